$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Noticias")

# Fix the image path text in E3
$ws.Range("E3").Value = "images/Peanut-Candy-assassinato/suspeito-crime-correção.png"

# Copy A2's/D2's original style (font2 rgb-black + border1) onto A3/D3
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# Copy A1's style (font1 theme + border1) onto A2 and A4
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)

# Copy D1's style onto D2 and D4
$ws.Range("D1").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 1 height adjustment (matches target diff)
$ws.Rows.Item(1).RowHeight = 19.5
